# Daily attendance processing - normalize "Recorded By" (column G) ordering.
# For every data row, the comma-separated list of recorders in column G is
# rotated by moving the last entry to the front (e.g. "A, B" -> "B, A"),
# except for the specific combination "System, backup@backdoor.com", which
# is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value()

    if ($null -eq $val -or $val -eq "") {
        continue
    }

    if ($val -eq "System, backup@backdoor.com") {
        continue
    }

    $parts = $val -split ", "
    if ($parts.Count -le 1) {
        continue
    }

    $rotated = (@($parts[-1]) + $parts[0..($parts.Count - 2)]) -join ", "

    if ($rotated -ne $val) {
        $cell.Value = $rotated
    }
}
